$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(13, 6).Value = 6423
$ws.Cells.Item(14, 6).Value = 1105
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(19, 6).Value = 283
$ws.Cells.Item(21, 6).Value = 215
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(23, 6).Value = 10352
$ws.Cells.Item(25, 6).Value = 10
$ws.Cells.Item(26, 6).Value = 1965
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 2225
$ws.Cells.Item(32, 6).Value = 177
$ws.Cells.Item(33, 6).Value = 24
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 1428
$ws.Cells.Item(38, 6).Value = 5348
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(40, 6).Value = 1209
$ws.Cells.Item(41, 6).Value = 706
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(46, 6).Value = 1001
$ws.Cells.Item(47, 6).Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 17
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(12, 6).Value = 9
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(18, 6).Value = 912
$ws.Cells.Item(19, 6).Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 9191
$ws.Cells.Item(6, 6).Value = 7187
$ws.Cells.Item(7, 6).Value = 187
$ws.Cells.Item(8, 6).Value = 51
$ws.Cells.Item(11, 6).Value = 5599
$ws.Cells.Item(12, 6).Value = 5599
$ws.Cells.Item(14, 6).Value = 6423
$ws.Cells.Item(15, 6).Value = 6423
$ws.Cells.Item(16, 6).Value = 1105
$ws.Cells.Item(17, 6).Value = 436
$ws.Cells.Item(18, 6).Value = 419
$ws.Cells.Item(19, 6).Value = 615
$ws.Cells.Item(21, 6).Value = 283
$ws.Cells.Item(22, 6).Value = 154
$ws.Cells.Item(23, 6).Value = 215
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 10352
$ws.Cells.Item(28, 6).Value = 1965
$ws.Cells.Item(30, 6).Value = 47
$ws.Cells.Item(31, 6).Value = 2225
$ws.Cells.Item(32, 6).Value = 85
$ws.Cells.Item(34, 6).Value = 177
$ws.Cells.Item(35, 6).Value = 24
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(38, 6).Value = 1428
$ws.Cells.Item(40, 6).Value = 5348
$ws.Cells.Item(41, 6).Value = 432
$ws.Cells.Item(42, 6).Value = 1209
$ws.Cells.Item(43, 6).Value = 706
$ws.Cells.Item(44, 6).Value = 127
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(48, 6).Value = 1001
$ws.Cells.Item(49, 6).Value = 1398
$ws.Cells.Item(50, 6).Value = 68
